# Update "想去人数" (want-to-go count) figures in column F across the
# "展览" (Exhibitions), "本地生活" (Local Life) and "全部类型" (All Types)
# sheets, matching the refreshed scrape captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 415
$wsExhibit.Range("F8").Value = 1865
$wsExhibit.Range("F9").Value = 801
$wsExhibit.Range("F10").Value = 20
$wsExhibit.Range("F12").Value = 1566
$wsExhibit.Range("F13").Value = 1566
$wsExhibit.Range("F14").Value = 1300
$wsExhibit.Range("F16").Value = 1373
$wsExhibit.Range("F18").Value = 392
$wsExhibit.Range("F23").Value = 7325
$wsExhibit.Range("F24").Value = 20
$wsExhibit.Range("F28").Value = 0
$wsExhibit.Range("F35").Value = 266
$wsExhibit.Range("F36").Value = 645
$wsExhibit.Range("F38").Value = 1345

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 236

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 236
$wsAll.Range("F11").Value = 415
$wsAll.Range("F13").Value = 1865
$wsAll.Range("F14").Value = 801
$wsAll.Range("F15").Value = 20
$wsAll.Range("F17").Value = 1566
$wsAll.Range("F18").Value = 1566
$wsAll.Range("F19").Value = 1300
$wsAll.Range("F21").Value = 1373
$wsAll.Range("F23").Value = 392
$wsAll.Range("F29").Value = 7325
$wsAll.Range("F31").Value = 228
$wsAll.Range("F35").Value = 266
$wsAll.Range("F38").Value = 645
$wsAll.Range("F43").Value = 1345
